# Generate Report for Handback
# Adds a new handback row (a00a7228-422b-48f1-b114-67c1f80c027f) to the
# "Overview", "zh-cn" and "de-de" sheets/tables, mirroring the pattern
# already used for the two existing rows.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overview" (table "Overview") -> new row 4
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A4").Value = "a00a7228-422b-48f1-b114-67c1f80c027f.md"
$wsOverview.Range("B4").Value = "e2e\a00a7228-422b-48f1-b114-67c1f80c027f.md"
$wsOverview.Range("C4").Value = ".md"
$wsOverview.Range("E4").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F4").Value = "Handed back: in sync with en-US"
$wsOverview.Range("G4").Value = "2016-10-19 16:49:16"
$wsOverview.Range("G4").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wb.Hyperlinks.Add($wsOverview.Range("B4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a00a7228422b48f1b11467c1f80c027fa00a722/e2e/a00a7228-422b-48f1-b114-67c1f80c027f.md", "", "", "e2e\a00a7228-422b-48f1-b114-67c1f80c027f.md")

$tblOverview = $wsOverview.ListObjects.Item(1)
$tblOverview.Resize($wsOverview.Range("A1:G4"))

# ---------------------------------------------------------------------
# Sheet "zh-cn" (table "zh-cn") -> new row 4
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A4").Value = "a00a7228-422b-48f1-b114-67c1f80c027f.md"
$wsZhCn.Range("B4").Value = ".md"
$wsZhCn.Range("C4").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("D4").Value = "e2e"
$wsZhCn.Range("E4").Value = "ht"
$wsZhCn.Range("F4").Value = "'True"
$wsZhCn.Range("G4").Value = "a00a7228-422b-48f1-b114-67c1f80c027f.73f25297e0291b637cf1e62f058f73fb1aaa5130.zh-cn.xlf"
$wsZhCn.Range("H4").Value = "2016-10-19 16:49:05"
$wsZhCn.Range("H4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("I4").Value = "a00a7228-422b-48f1-b114-67c1f80c027f.md"
$wsZhCn.Range("J4").Value = "a00a7228-422b-48f1-b114-67c1f80c027f.73f25297e0291b637cf1e62f058f73fb1aaa5130.zh-cn.xlf"
$wsZhCn.Range("K4").Value = "2016-10-19 16:49:47"
$wsZhCn.Range("K4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("L4").Value = "'"
$wsZhCn.Range("M4").Value = "'True"
$wsZhCn.Range("N4").Value = "'"
$wsZhCn.Range("O4").Value = "'False"
$wsZhCn.Range("P4").Value = "'"

$wb.Hyperlinks.Add($wsZhCn.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a00a7228422b48f1b11467c1f80c027fa00a722/e2e/a00a7228-422b-48f1-b114-67c1f80c027f.md", "", "", "a00a7228-422b-48f1-b114-67c1f80c027f.md")
$wb.Hyperlinks.Add($wsZhCn.Range("I4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/a00a7228422b48f1b11467c1f80c027fa00a722/e2e/a00a7228-422b-48f1-b114-67c1f80c027f.md", "", "", "a00a7228-422b-48f1-b114-67c1f80c027f.md")

$tblZhCn = $wsZhCn.ListObjects.Item(1)
$tblZhCn.Resize($wsZhCn.Range("A1:P4"))

# ---------------------------------------------------------------------
# Sheet "de-de" (table "de-de") -> new row 4
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A4").Value = "a00a7228-422b-48f1-b114-67c1f80c027f.md"
$wsDeDe.Range("B4").Value = ".md"
$wsDeDe.Range("C4").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("D4").Value = "e2e"
$wsDeDe.Range("E4").Value = "ht"
$wsDeDe.Range("F4").Value = "'True"
$wsDeDe.Range("G4").Value = "a00a7228-422b-48f1-b114-67c1f80c027f.73f25297e0291b637cf1e62f058f73fb1aaa5130.de-de.xlf"
$wsDeDe.Range("H4").Value = "2016-10-19 16:49:16"
$wsDeDe.Range("H4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("I4").Value = "a00a7228-422b-48f1-b114-67c1f80c027f.md"
$wsDeDe.Range("J4").Value = "a00a7228-422b-48f1-b114-67c1f80c027f.73f25297e0291b637cf1e62f058f73fb1aaa5130.de-de.xlf"
$wsDeDe.Range("K4").Value = "2016-10-19 16:50:06"
$wsDeDe.Range("K4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("L4").Value = "'"
$wsDeDe.Range("M4").Value = "'True"
$wsDeDe.Range("N4").Value = "'"
$wsDeDe.Range("O4").Value = "'False"
$wsDeDe.Range("P4").Value = "'"

$wb.Hyperlinks.Add($wsDeDe.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a00a7228422b48f1b11467c1f80c027fa00a722/e2e/a00a7228-422b-48f1-b114-67c1f80c027f.md", "", "", "a00a7228-422b-48f1-b114-67c1f80c027f.md")
$wb.Hyperlinks.Add($wsDeDe.Range("I4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/a00a7228422b48f1b11467c1f80c027fa00a722/e2e/a00a7228-422b-48f1-b114-67c1f80c027f.md", "", "", "a00a7228-422b-48f1-b114-67c1f80c027f.md")

$tblDeDe = $wsDeDe.ListObjects.Item(1)
$tblDeDe.Resize($wsDeDe.Range("A1:P4"))

"Report row added for a00a7228-422b-48f1-b114-67c1f80c027f.md"
